$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing E-column values (rows 2-61) ---
$ws.Range("E5").Value = 20
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 7
$ws.Range("E14").Value = 5
$ws.Range("E15").Value = 2
$ws.Range("E16").Value = 6
$ws.Range("E17").Value = 2
$ws.Range("E18").Value = 6
$ws.Range("E19").Value = 9
$ws.Range("E26").Value = 3
$ws.Range("E27").Value = 8
$ws.Range("E28").Value = 1
$ws.Range("E29").Value = 6
$ws.Range("E30").Value = 2
$ws.Range("E31").Value = 6
$ws.Range("E38").Value = 10
$ws.Range("E39").Value = 4
$ws.Range("E40").Value = 5
$ws.Range("E41").Value = 4
$ws.Range("E42").Value = 3
$ws.Range("E43").Value = 9
$ws.Range("E50").Value = 3
$ws.Range("E51").Value = 5
$ws.Range("E52").Value = 10
$ws.Range("E53").Value = 10
$ws.Range("E54").Value = 1
$ws.Range("E55").Value = 2

# --- Append new rows 62-109 ---
$ws.Range("A62").Value = "C1"
$ws.Range("B62").Value = "Male"
$ws.Range("C62").Value = "20-40"
$ws.Range("D62").Value = 2024
$ws.Range("E62").Value = 3
$ws.Range("A63").Value = "C1"
$ws.Range("B63").Value = "Male"
$ws.Range("C63").Value = "40-60"
$ws.Range("D63").Value = 2024
$ws.Range("E63").Value = 5
$ws.Range("A64").Value = "C1"
$ws.Range("B64").Value = "Male"
$ws.Range("C64").Value = "60-80"
$ws.Range("D64").Value = 2024
$ws.Range("E64").Value = 10
$ws.Range("A65").Value = "C1"
$ws.Range("B65").Value = "Female"
$ws.Range("C65").Value = "20-40"
$ws.Range("D65").Value = 2024
$ws.Range("E65").Value = 10
$ws.Range("A66").Value = "C1"
$ws.Range("B66").Value = "Female"
$ws.Range("C66").Value = "40-60"
$ws.Range("D66").Value = 2024
$ws.Range("E66").Value = 1
$ws.Range("A67").Value = "C1"
$ws.Range("B67").Value = "Female"
$ws.Range("C67").Value = "60-80"
$ws.Range("D67").Value = 2024
$ws.Range("E67").Value = 2
$ws.Range("A68").Value = "C2"
$ws.Range("B68").Value = "Male"
$ws.Range("C68").Value = "20-40"
$ws.Range("D68").Value = 2024
$ws.Range("E68").Value = 7
$ws.Range("A69").Value = "C2"
$ws.Range("B69").Value = "Male"
$ws.Range("C69").Value = "40-60"
$ws.Range("D69").Value = 2024
$ws.Range("E69").Value = 0
$ws.Range("A70").Value = "C2"
$ws.Range("B70").Value = "Male"
$ws.Range("C70").Value = "60-80"
$ws.Range("D70").Value = 2024
$ws.Range("E70").Value = 4
$ws.Range("A71").Value = "C2"
$ws.Range("B71").Value = "Female"
$ws.Range("C71").Value = "20-40"
$ws.Range("D71").Value = 2024
$ws.Range("E71").Value = 0
$ws.Range("A72").Value = "C2"
$ws.Range("B72").Value = "Female"
$ws.Range("C72").Value = "40-60"
$ws.Range("D72").Value = 2024
$ws.Range("E72").Value = 8
$ws.Range("A73").Value = "C2"
$ws.Range("B73").Value = "Female"
$ws.Range("C73").Value = "60-80"
$ws.Range("D73").Value = 2024
$ws.Range("E73").Value = 7
$ws.Range("A74").Value = "C1"
$ws.Range("B74").Value = "Male"
$ws.Range("C74").Value = "20-40"
$ws.Range("D74").Value = 2025
$ws.Range("E74").Value = 3
$ws.Range("A75").Value = "C1"
$ws.Range("B75").Value = "Male"
$ws.Range("C75").Value = "40-60"
$ws.Range("D75").Value = 2025
$ws.Range("E75").Value = 5
$ws.Range("A76").Value = "C1"
$ws.Range("B76").Value = "Male"
$ws.Range("C76").Value = "60-80"
$ws.Range("D76").Value = 2025
$ws.Range("E76").Value = 10
$ws.Range("A77").Value = "C1"
$ws.Range("B77").Value = "Female"
$ws.Range("C77").Value = "20-40"
$ws.Range("D77").Value = 2025
$ws.Range("E77").Value = 10
$ws.Range("A78").Value = "C1"
$ws.Range("B78").Value = "Female"
$ws.Range("C78").Value = "40-60"
$ws.Range("D78").Value = 2025
$ws.Range("E78").Value = 1
$ws.Range("A79").Value = "C1"
$ws.Range("B79").Value = "Female"
$ws.Range("C79").Value = "60-80"
$ws.Range("D79").Value = 2025
$ws.Range("E79").Value = 2
$ws.Range("A80").Value = "C2"
$ws.Range("B80").Value = "Male"
$ws.Range("C80").Value = "20-40"
$ws.Range("D80").Value = 2025
$ws.Range("E80").Value = 7
$ws.Range("A81").Value = "C2"
$ws.Range("B81").Value = "Male"
$ws.Range("C81").Value = "40-60"
$ws.Range("D81").Value = 2025
$ws.Range("E81").Value = 0
$ws.Range("A82").Value = "C2"
$ws.Range("B82").Value = "Male"
$ws.Range("C82").Value = "60-80"
$ws.Range("D82").Value = 2025
$ws.Range("E82").Value = 4
$ws.Range("A83").Value = "C2"
$ws.Range("B83").Value = "Female"
$ws.Range("C83").Value = "20-40"
$ws.Range("D83").Value = 2025
$ws.Range("E83").Value = 0
$ws.Range("A84").Value = "C2"
$ws.Range("B84").Value = "Female"
$ws.Range("C84").Value = "40-60"
$ws.Range("D84").Value = 2025
$ws.Range("E84").Value = 8
$ws.Range("A85").Value = "C2"
$ws.Range("B85").Value = "Female"
$ws.Range("C85").Value = "60-80"
$ws.Range("D85").Value = 2025
$ws.Range("E85").Value = 7
$ws.Range("A86").Value = "C1"
$ws.Range("B86").Value = "Male"
$ws.Range("C86").Value = "20-40"
$ws.Range("E86").Value = 3
$ws.Range("A87").Value = "C1"
$ws.Range("B87").Value = "Male"
$ws.Range("C87").Value = "40-60"
$ws.Range("E87").Value = 5
$ws.Range("A88").Value = "C1"
$ws.Range("B88").Value = "Male"
$ws.Range("C88").Value = "60-80"
$ws.Range("E88").Value = 10
$ws.Range("A89").Value = "C1"
$ws.Range("B89").Value = "Female"
$ws.Range("C89").Value = "20-40"
$ws.Range("E89").Value = 10
$ws.Range("A90").Value = "C1"
$ws.Range("B90").Value = "Female"
$ws.Range("C90").Value = "40-60"
$ws.Range("E90").Value = 1
$ws.Range("A91").Value = "C1"
$ws.Range("B91").Value = "Female"
$ws.Range("C91").Value = "60-80"
$ws.Range("E91").Value = 2
$ws.Range("A92").Value = "C2"
$ws.Range("B92").Value = "Male"
$ws.Range("C92").Value = "20-40"
$ws.Range("E92").Value = 7
$ws.Range("A93").Value = "C2"
$ws.Range("B93").Value = "Male"
$ws.Range("C93").Value = "40-60"
$ws.Range("E93").Value = 0
$ws.Range("A94").Value = "C2"
$ws.Range("B94").Value = "Male"
$ws.Range("C94").Value = "60-80"
$ws.Range("E94").Value = 4
$ws.Range("A95").Value = "C2"
$ws.Range("B95").Value = "Female"
$ws.Range("C95").Value = "20-40"
$ws.Range("E95").Value = 0
$ws.Range("A96").Value = "C2"
$ws.Range("B96").Value = "Female"
$ws.Range("C96").Value = "40-60"
$ws.Range("E96").Value = 8
$ws.Range("A97").Value = "C2"
$ws.Range("B97").Value = "Female"
$ws.Range("C97").Value = "60-80"
$ws.Range("E97").Value = 7
$ws.Range("A98").Value = "C1"
$ws.Range("B98").Value = "Male"
$ws.Range("C98").Value = "20-40"
$ws.Range("D98").Value = 2026
$ws.Range("E98").Value = 3
$ws.Range("A99").Value = "C1"
$ws.Range("B99").Value = "Male"
$ws.Range("C99").Value = "40-60"
$ws.Range("D99").Value = 2026
$ws.Range("E99").Value = 5
$ws.Range("A100").Value = "C1"
$ws.Range("B100").Value = "Male"
$ws.Range("C100").Value = "60-80"
$ws.Range("D100").Value = 2026
$ws.Range("E100").Value = 10
$ws.Range("A101").Value = "C1"
$ws.Range("B101").Value = "Female"
$ws.Range("C101").Value = "20-40"
$ws.Range("D101").Value = 2026
$ws.Range("E101").Value = 10
$ws.Range("A102").Value = "C1"
$ws.Range("B102").Value = "Female"
$ws.Range("C102").Value = "40-60"
$ws.Range("D102").Value = 2026
$ws.Range("E102").Value = 1
$ws.Range("A103").Value = "C1"
$ws.Range("B103").Value = "Female"
$ws.Range("C103").Value = "60-80"
$ws.Range("D103").Value = 2026
$ws.Range("E103").Value = 2
$ws.Range("A104").Value = "C2"
$ws.Range("B104").Value = "Male"
$ws.Range("C104").Value = "20-40"
$ws.Range("D104").Value = 2026
$ws.Range("E104").Value = 7
$ws.Range("A105").Value = "C2"
$ws.Range("B105").Value = "Male"
$ws.Range("C105").Value = "40-60"
$ws.Range("D105").Value = 2026
$ws.Range("E105").Value = 0
$ws.Range("A106").Value = "C2"
$ws.Range("B106").Value = "Male"
$ws.Range("C106").Value = "60-80"
$ws.Range("D106").Value = 2026
$ws.Range("E106").Value = 4
$ws.Range("A107").Value = "C2"
$ws.Range("B107").Value = "Female"
$ws.Range("C107").Value = "20-40"
$ws.Range("D107").Value = 2026
$ws.Range("E107").Value = 0
$ws.Range("A108").Value = "C2"
$ws.Range("B108").Value = "Female"
$ws.Range("C108").Value = "40-60"
$ws.Range("D108").Value = 2026
$ws.Range("E108").Value = 8
$ws.Range("A109").Value = "C2"
$ws.Range("B109").Value = "Female"
$ws.Range("C109").Value = "60-80"
$ws.Range("D109").Value = 2026
$ws.Range("E109").Value = 7

Write-Host "Edit complete"
